$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [new C value, new D value]
$updates = @{
    2  = @(61, 1.05)
    3  = @(90, 1.55)
    5  = @(48, 0.83)
    6  = @(254, 4.38)
    7  = @(169, 2.91)
    8  = @(1153, 19.86)
    9  = @(172, 2.96)
    11 = @(65, 1.12)
    12 = @(379, 6.53)
    13 = @(146, 2.52)
    14 = @(67, 1.15)
    15 = @(223, 3.84)
    16 = @(94, 1.62)
    17 = @(269, 4.63)
    18 = @(209, 3.6)
    20 = @(588, 10.13)
    21 = @(85, 1.46)
    22 = @(45, 0.78)
    23 = @(58, 1)
    24 = @(180, 3.1)
    25 = @(271, 4.67)
    27 = @(332, 5.72)
    28 = @(59, 1.02)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
